$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 7143132.5
$ws.Range("I33").Value = 7143132.5
$ws.Range("K33").Value = 7143132.5
$ws.Range("M33").Value = -7142903.5

# Row 40
$ws.Range("H40").Value = 1379.4897
$ws.Range("I40").Value = 1188.8334
$ws.Range("K40").Value = 1188.8334
$ws.Range("M40").Value = -1013.8334

# Row 64
$ws.Range("H64").Value = 8439.4
$ws.Range("I64").Value = 9049.25
$ws.Range("J64").Value = 6000
$ws.Range("K64").Value = 9049.25
$ws.Range("L64").Value = 6000
$ws.Range("M64").Value = -8801.25
$ws.Range("N64").Value = -6496

# Row 67
$ws.Range("H67").Value = 8439.4
$ws.Range("I67").Value = 9049.25
$ws.Range("J67").Value = 6000
$ws.Range("K67").Value = 9049.25
$ws.Range("L67").Value = 6000
$ws.Range("M67").Value = -8191.25
$ws.Range("N67").Value = -7716

# Row 96
$ws.Range("H96").Value = 703.1111
$ws.Range("J96").Value = 818.5
$ws.Range("L96").Value = 2455.5
$ws.Range("N96").Value = -5201.5

# Row 98
$ws.Range("H98").Value = 699.8421
$ws.Range("I98").Value = 749.4286
$ws.Range("J98").Value = 561
$ws.Range("K98").Value = 749.4286
$ws.Range("L98").Value = 561
$ws.Range("M98").Value = 748.5714
$ws.Range("N98").Value = -3557

# Row 103
$ws.Range("H103").Value = 862.25
$ws.Range("I103").Value = 500
$ws.Range("J103").Value = 983
$ws.Range("K103").Value = 1500
$ws.Range("L103").Value = 2949
$ws.Range("M103").Value = -914
$ws.Range("N103").Value = -4121

# Row 106
$ws.Range("H106").Value = 1625
$ws.Range("I106").Value = 1032.1428
$ws.Range("K106").Value = 1032.1428
$ws.Range("M106").Value = -401.1428000000001

# Row 116
$ws.Range("H116").Value = 8276
$ws.Range("J116").Value = 5489.3335
$ws.Range("L116").Value = 5489.3335
$ws.Range("N116").Value = -12373.3335

# Row 122
$ws.Range("H122").Value = 699.8421
$ws.Range("I122").Value = 749.4286
$ws.Range("J122").Value = 561
$ws.Range("K122").Value = 2248.2858
$ws.Range("L122").Value = 1683
$ws.Range("M122").Value = 201.7142000000003
$ws.Range("N122").Value = -6583

# Row 132
$ws.Range("H132").Value = 6192308
$ws.Range("I132").Value = 7599389.5
$ws.Range("K132").Value = 22798168.5
$ws.Range("M132").Value = -22795638.5

# Row 137
$ws.Range("H137").Value = 1959
$ws.Range("I137").Value = 2523.5
$ws.Range("J137").Value = 1313.8572
$ws.Range("K137").Value = 7570.5
$ws.Range("L137").Value = 3941.5716
$ws.Range("M137").Value = -5020.5
$ws.Range("N137").Value = -9041.571599999999


# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 6913.7407
$ws.Range("I86").Value = 5283.2
$ws.Range("K86").Value = 5283.2
$ws.Range("M86").Value = -4160.2

# Row 89
$ws.Range("H89").Value = 6913.7407
$ws.Range("I89").Value = 5283.2
$ws.Range("K89").Value = 26416
$ws.Range("M89").Value = -20800

# Row 107
$ws.Range("H107").Value = 1790.1951
$ws.Range("I107").Value = 1855.8485
$ws.Range("J107").Value = 1519.375
$ws.Range("K107").Value = 1855.8485
$ws.Range("L107").Value = 1519.375
$ws.Range("M107").Value = 64.15149999999994
$ws.Range("N107").Value = -5359.375

# Row 132
$ws.Range("H132").Value = 80497.5
$ws.Range("J132").Value = 80497.5
$ws.Range("L132").Value = 80497.5
$ws.Range("N132").Value = -90617.5


# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 3077
$ws.Range("I16").Value = 4818.6
$ws.Range("K16").Value = 4818.6
$ws.Range("M16").Value = -4531.6

# Row 22
$ws.Range("H22").Value = 709.8
$ws.Range("J22").Value = 983
$ws.Range("L22").Value = 983
$ws.Range("N22").Value = -1683

# Row 31
$ws.Range("H31").Value = 77435.28999999999
$ws.Range("I31").Value = 101972.5
$ws.Range("J31").Value = 16092.25
$ws.Range("K31").Value = 101972.5
$ws.Range("L31").Value = 16092.25
$ws.Range("M31").Value = -101677.5
$ws.Range("N31").Value = -16682.25

# Row 34
$ws.Range("H34").Value = 77435.28999999999
$ws.Range("I34").Value = 101972.5
$ws.Range("J34").Value = 16092.25
$ws.Range("K34").Value = 101972.5
$ws.Range("L34").Value = 16092.25
$ws.Range("M34").Value = -101770.5
$ws.Range("N34").Value = -16496.25

# Row 58
$ws.Range("H58").Value = 3400.5
$ws.Range("I58").Value = 3479.8
$ws.Range("J58").Value = 3004
$ws.Range("K58").Value = 3479.8
$ws.Range("L58").Value = 3004
$ws.Range("M58").Value = -3276.8
$ws.Range("N58").Value = -3410

# Row 94
$ws.Range("H94").Value = 1163.2632
$ws.Range("I94").Value = 1305.6154
$ws.Range("K94").Value = 1305.6154
$ws.Range("M94").Value = -854.6153999999999

# Row 113
$ws.Range("H113").Value = 3077
$ws.Range("I113").Value = 4818.6
$ws.Range("K113").Value = 4818.6
$ws.Range("M113").Value = -2648.6

# Row 122
$ws.Range("H122").Value = 1110.1818
$ws.Range("I122").Value = 1055.6666
$ws.Range("K122").Value = 3166.9998
$ws.Range("M122").Value = -716.9998000000001

# Row 134
$ws.Range("H134").Value = 15687.657
$ws.Range("I134").Value = 9166.75
$ws.Range("K134").Value = 27500.25
$ws.Range("M134").Value = -24965.25

# Row 136
$ws.Range("H136").Value = 3400.5
$ws.Range("I136").Value = 3479.8
$ws.Range("J136").Value = 3004
$ws.Range("K136").Value = 10439.4
$ws.Range("L136").Value = 9012
$ws.Range("M136").Value = -7889.400000000001
$ws.Range("N136").Value = -14112


# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 7144485.5
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 7144485.5
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 21433456.5
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -21435078.5

# Row 71
$ws.Range("H71").Value = 7144485.5
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 7144485.5
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 64300369.5
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -64308481.5


# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 34
$ws.Range("H34").Value = 37500
$ws.Range("J34").Value = 37500
$ws.Range("L34").Value = 37500
$ws.Range("N34").Value = -38036

# Row 70
$ws.Range("H70").Value = 15713.667
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 15713.667
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 15713.667
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -16253.667

# Row 73
$ws.Range("H73").Value = 15713.667
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 15713.667
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 15713.667
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -17585.667

# Row 76
$ws.Range("H76").Value = 37500
$ws.Range("J76").Value = 37500
$ws.Range("L76").Value = 37500
$ws.Range("N76").Value = -38130

# Row 79
$ws.Range("H79").Value = 37500
$ws.Range("J79").Value = 37500
$ws.Range("L79").Value = 37500
$ws.Range("N79").Value = -39684

# Row 122
$ws.Range("H122").Value = 2952.4707
$ws.Range("I122").Value = 2520.4
$ws.Range("K122").Value = 7561.200000000001
$ws.Range("M122").Value = -5111.200000000001

# Row 123
$ws.Range("H123").Value = 26060
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 26060
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 26060
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -30960


# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 22733536
$ws.Range("I16").Value = 100002450
$ws.Range("J16").Value = 7385.5884
$ws.Range("K16").Value = 100002450
$ws.Range("L16").Value = 7385.5884
$ws.Range("M16").Value = -100002280
$ws.Range("N16").Value = -7725.5884

# Row 40
$ws.Range("H40").Value = 3889.6956
$ws.Range("I40").Value = 3627.5
$ws.Range("J40").Value = 4489
$ws.Range("K40").Value = 3627.5
$ws.Range("L40").Value = 4489
$ws.Range("M40").Value = -3491.5
$ws.Range("N40").Value = -4761

# Row 100
$ws.Range("H100").Value = 440551.12
$ws.Range("I100").Value = 530451.4
$ws.Range("J100").Value = 13525
$ws.Range("K100").Value = 530451.4
$ws.Range("L100").Value = 13525
$ws.Range("M100").Value = -529910.4
$ws.Range("N100").Value = -14607

# Row 122
$ws.Range("H122").Value = 4920.8125
$ws.Range("I122").Value = 4457.4546
$ws.Range("J122").Value = 5940.2
$ws.Range("K122").Value = 13372.3638
$ws.Range("L122").Value = 17820.6
$ws.Range("M122").Value = -10922.3638
$ws.Range("N122").Value = -22720.6

# Row 132
$ws.Range("H132").Value = 4859.591
$ws.Range("I132").Value = 3743.3572
$ws.Range("J132").Value = 6813
$ws.Range("K132").Value = 11230.0716
$ws.Range("L132").Value = 20439
$ws.Range("M132").Value = -8700.071599999999
$ws.Range("N132").Value = -25499


# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 3768.9565
$ws.Range("I132").Value = 3769.3
$ws.Range("K132").Value = 11307.9
$ws.Range("M132").Value = -8777.900000000001

# Row 136
$ws.Range("H136").Value = 2179
$ws.Range("I136").Value = 2180.8823
$ws.Range("J136").Value = 2175
$ws.Range("K136").Value = 6542.646900000001
$ws.Range("L136").Value = 6525
$ws.Range("M136").Value = -3992.646900000001
$ws.Range("N136").Value = -11625

